$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying "expression" values in column B were re-labelled to a
# single "neutral" class across the whole existing dataset (rows 2-35),
# and the dataset was extended with additional rows (36-68) continuing
# the index sequence in column A, all likewise labelled "neutral".

# 1) Extend column A (the numeric index) down to row 68, continuing the
#    existing 0-based sequence (row 36 -> 34 ... row 68 -> 66).
#    First clone the formatting already used by the index column so the
#    new cells match the existing look (centered, bordered, bold).
$ws.Range("A35").Copy()
$ws.Range("A36:A68").PasteSpecial(-4122)

for ($r = 36; $r -le 68; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# 2) Normalize every label in column B (existing rows plus the newly
#    added ones) to "neutral".
$ws.Range("B2:B68").Value = "neutral"
